$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts B:E -> C:F)
$ws.Columns("B").Insert()

# New header cell for the inserted "Internal" column
$ws.Range("B2").Value = "Internal"

# Second header block (row 4) also gets the "Internal" label
$ws.Range("B4").Value = "Internal"

# Data rows for the new Internal column
$ws.Range("B3").Value = "FALSE"
$ws.Range("B5").Value = "FALSE"
$ws.Range("B6").Value = "FALSE"

# Update the selection to match the post-edit state
$ws.Range("B4:B6").Select()
